$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @("LP-043273","LP-045293","LP-045321","LP-045512","LP-045682","LP-045825","LP-046160","LP-046182","LP-046995","LP-047091","LP-047093","LP-047094","LP-047112","LP-047246","LP-047250","LP-047330","LP-047339","LP-047478","LP-047480","LP-047481","LP-047482","LP-047483","LP-047484","LP-047485","LP-047486","LP-047487","LP-047503","LP-047521","LP-047505","LP-047564","LP-047770","LP-047774","LP-047833","LP-047874","LP-047876","LP-047877","LP-047891","LP-047892","LP-048036","LP-048151","LP-048271","LP-048275","LP-048276","LP-048292","LP-048368","LP-048382","LP-048589","LP-048597","LP-048653","LP-048654","LP-048704")

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $values[$i]
}
